# Update "想去人数" (F column) counts on three sheets to reflect newly
# regenerated site data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2794
$ws1.Range("F20").Value = 505
$ws1.Range("F22").Value = 1278
$ws1.Range("F25").Value = 2039
$ws1.Range("F26").Value = 156
$ws1.Range("F29").Value = 3232

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 1464
$ws3.Range("F12").Value = 606

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1464
$ws4.Range("F14").Value = 2794
$ws4.Range("F17").Value = 606
$ws4.Range("F30").Value = 505
$ws4.Range("F38").Value = 1278
$ws4.Range("F43").Value = 2039
$ws4.Range("F46").Value = 156
$ws4.Range("F49").Value = 3232
